# Daily update at 8 AM UTC
# Append the next day's data as a new row at the bottom of the data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = 46001
$ws.Cells.Item($newRow, 2).Value = 117
$ws.Cells.Item($newRow, 3).Value = 127
$ws.Cells.Item($newRow, 4).Value = 116

# Carry the date-formatted style from the previous row's date cell.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
